$wb = $excel.ActiveWorkbook

# Insert a new worksheet named "Contact Type" right after the "Contact" sheet
# and before "Industry Group" (matches sheetId=6, rId3, pushing Industry Group to rId4).
$contactSheet = $wb.Worksheets.Item("Contact")
$newSheet = $wb.Worksheets.Add([System.Type]::Missing, $contactSheet)
$newSheet.Name = "Contact Type"

# Populate the new sheet with its header + single data row.
$newSheet.Range("A1").Value = "Contact Type"
$newSheet.Range("A2").Value = "Houlihan Employee"

# Bold header cell to match the style used on the other sheets' header rows.
$newSheet.Range("A1").Font.Bold = $true

# Match the column width used in the authored workbook (stored width 18).
$newSheet.Columns.Item(1).ColumnWidth = 17.1666666666667

# Update the selection on the "Contact" sheet: the whole column B is now selected
# and it is no longer the active/tab-selected sheet.
$contactSheet.Range("B1:B1048576").Select() | Out-Null

# Make "Contact Type" the active sheet and set its selection to B8,
# which also marks it as the tab-selected sheet (activeTab becomes index 2).
$newSheet.Activate()
$newSheet.Range("B8").Select() | Out-Null
